# Add a new weekly price record for "Ají" (Cristal variety) at Terminal
# Hortofrutícola Agro Chillán. Insert a new row at row 22 (shifting the
# existing rows 22:82 down to 23:83) and populate it with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("22:22").Insert()

$ws.Range("A22").Value = 7
$ws.Range("B22").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C22").Value = "Ñuble"
$ws.Range("D22").Value = 44672
$ws.Range("E22").Value = 16
$ws.Range("F22").Value = 100112021
$ws.Range("G22").Value = "Ají"
$ws.Range("H22").Value = "Cristal"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 60
$ws.Range("K22").Value = 15500
$ws.Range("L22").Value = 16000
$ws.Range("M22").Value = 15750
$ws.Range("N22").Value = "$/saco 25 kilos"
$ws.Range("O22").Value = "Región del Maule"
$ws.Range("P22").Value = 630
$ws.Range("Q22").Value = 25
$ws.Range("R22").Value = "Hortaliza"
